$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Footer/date placeholders on the slide master and every slide layout:
#    "3/27/2020" -> "3/28/2020"
# ---------------------------------------------------------------------------
$master = $p.Slides.Item(1).Master

$masterDate = $master.Shapes.Item(4)
if ($masterDate.HasTextFrame -and $masterDate.TextFrame.TextRange.Text -eq "3/27/2020") {
    $masterDate.TextFrame.TextRange.Text = "3/28/2020"
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    $dateShape = $layout.Shapes.Item(4)
    if ($dateShape.HasTextFrame -and $dateShape.TextFrame.TextRange.Text -eq "3/27/2020") {
        $dateShape.TextFrame.TextRange.Text = "3/28/2020"
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 3 ("Concept"): fix wording in the description paragraph
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(2)
$tr3 = $shp3.TextFrame.TextRange
$full3 = $tr3.Text
$old3 = "`tOur site allows clients to put items into an online shopping cart and then proceed to a `tcheckout with those items."
$new3 = "`tOur application allows clients to put items into an online shopping cart and then proceed  `tto a checkout with those items."
$idx3 = $full3.IndexOf($old3)
if ($idx3 -ge 0) {
    $sub3 = $tr3.Characters($idx3 + 1, $old3.Length)
    $sub3.Text = $new3
}

# ---------------------------------------------------------------------------
# 3) Slide 4 ("Process"): fix typo + merge the three runs of the last bullet
#    into a single run
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
$tr4 = $shp4.TextFrame.TextRange
$full4 = $tr4.Text
$old4 = "`tFigured out hot to put admin permissions on certain pages of the website and how to change `tthe website based on who is logged in."
$new4 = "`tFigured out how to put admin permissions on certain pages of the website and how to `tchange `tthe website based on who is logged in."
$idx4 = $full4.IndexOf($old4)
if ($idx4 -ge 0) {
    $sub4 = $tr4.Characters($idx4 + 1, $old4.Length)
    $sub4.Text = $new4
}
